# Adds three new daily NRFI-tracker sheets ("08-22-24", "08-23-24", "08-29-24")
# to the end of the workbook, matching the layout/format used by the existing
# "Games"/"Score" day sheets.

$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{
        Name = "08-22-24"
        Rows = @(
            @("('NYM', 'SD')", 0.728),
            @("('MIL', 'STL')", 0.719),
            @("('COL', 'WSH')", 0.713),
            @("('CIN', 'PIT')", 0.7),
            @("('CHC', 'DET')", 0.693),
            @("('OAK', 'TB')", 0.61),
            @("('BAL', 'HOU')", 0.572),
            @("('ATL', 'PHI')", 0.556),
            @("('CLE', 'NYY')", 0.456),
            @("('LAA', 'TOR')", -0.022)
        )
    },
    @{
        Name = "08-23-24"
        Rows = @(
            @("('CLE', 'TEX')", 0.725),
            @("('SEA', 'SF')", 0.719),
            @("('MIN', 'STL')", 0.717),
            @("('KC', 'PHI')", 0.707),
            @("('CHC', 'MIA')", 0.691),
            @("('MIL', 'OAK')", 0.68),
            @("('ATL', 'WSH')", 0.647),
            @("('CIN', 'PIT')", 0.635),
            @("('AZ', 'BOS')", 0.616),
            @("('COL', 'NYY')", 0.539),
            @("('NYM', 'SD')", 0.346),
            @("('CWS', 'DET')", 0.23),
            @("('LAD', 'TB')", 0.115),
            @("('BAL', 'HOU')", 0.003),
            @("('LAA', 'TOR')", -0.117)
        )
    },
    @{
        Name = "08-29-24"
        Rows = @(
            @("('COL', 'MIA')", 0.723),
            @("('CWS', 'TEX')", 0.714),
            @("('ATL', 'PHI')", 0.712),
            @("('HOU', 'KC')", 0.628),
            @("('BAL', 'LAD')", 0.421),
            @("('SD', 'STL')", 0.387),
            @("('CIN', 'OAK')", 0.33),
            @("('BOS', 'TOR')", 0.267),
            @("('MIL', 'SF')", 0.261),
            @("('AZ', 'NYM')", 0.252),
            @("('DET', 'LAA')", 0.051)
        )
    }
)

# Name of the sheet currently last in tab order -- new sheets get moved to
# just after this one, one at a time, so the three new sheets end up in
# order at the very end of the workbook.
$lastSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name

foreach ($sheetInfo in $sheetsData) {
    $ws = $wb.Worksheets.Add()
    $ws.Name = $sheetInfo.Name

    # Header row
    $ws.Range("A1").Value = "Games"
    $ws.Range("B1").Value = "Score"
    $header = $ws.Range("A1:B1")
    $header.Font.Bold = $true
    $header.HorizontalAlignment = -4108
    $header.VerticalAlignment = -4160
    $header.Borders.LineStyle = 1

    # Data rows
    $r = 2
    foreach ($row in $sheetInfo.Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $r = $r + 1
    }

    # Move the new sheet to the end of the workbook (right after the
    # previously-last sheet), then remember its name as the new tail.
    $target = $wb.Worksheets.Item($lastSheetName)
    $ws.Move($null, $target)
    $lastSheetName = $sheetInfo.Name
}
